$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New image URL column (G) for existing rows 2-4
$ws.Range("G2").Value = "https://images.stockx.com/images/adidas-Crazy-BYW-X-2-Ubiq.png?fit=fill&bg=FFFFFF&w=700&h=500&auto=format,compress&q=90&dpr=2&trim=color&updated_at=1606936652"
$ws.Range("G3").Value = "https://images.stockx.com/images/adidas-Pro-Bounce-Mid-2018-Kristaps-Porzingis-PE.png?fit=fill&bg=FFFFFF&w=700&h=500&auto=format,compress&q=90&dpr=2&trim=color&updated_at=1626898739"
$ws.Range("G4").Value = "https://images.stockx.com/images/Adidas-Ozweego-3-Raf-Simons-Cream-White-Core-Black-Product.jpg?fit=fill&bg=FFFFFF&w=700&h=500&auto=format,compress&q=90&dpr=2&trim=color&updated_at=1619145566"

# Row 6: adidas Poweralley 4
$ws.Range("B6").Value = "adidas Poweralley 4"
$ws.Range("C6").Value = "Core Black/White/Gold"
$ws.Range("D6").Value = "Not Found"
$ws.Range("E6").Value = "Not Found"
$ws.Range("F6").Value = "adidas"
$ws.Range("G6").Value = "https://images.stockx.com/images/adidas-Poweralley-4-Core-Black-White-Gold.jpg?fit=fill&bg=FFFFFF&w=700&h=500&auto=format,compress&q=90&dpr=2&trim=color&updated_at=1627414789"

# Row 7: adidas ZX 4000 (Price/Date look numeric - force text with leading apostrophe)
$ws.Range("B7").Value = "adidas ZX 4000"
$ws.Range("C7").Value = "Sesame/Clear Brown/Core Black"
$ws.Range("D7").Value = "'120"
$ws.Range("E7").Value = "'2019-03-08"
$ws.Range("F7").Value = "adidas"
$ws.Range("G7").Value = "https://images.stockx.com/images/adidas-ZX-4000-Sesame.png?fit=fill&bg=FFFFFF&w=700&h=500&auto=format,compress&q=90&dpr=2&trim=color&updated_at=1627414851"

# Row 8: adidas ZX 500 OG (Price looks numeric - force text with leading apostrophe)
$ws.Range("B8").Value = "adidas ZX 500 OG"
$ws.Range("C8").Value = "Grey/Red"
$ws.Range("D8").Value = "'140"
$ws.Range("E8").Value = "Not Found"
$ws.Range("F8").Value = "adidas"
$ws.Range("G8").Value = "https://images.stockx.com/images/adidas-ZX-500-OG-Grey-Red.jpg?fit=fill&bg=FFFFFF&w=700&h=500&auto=format,compress&q=90&dpr=2&trim=color&updated_at=1627415526"

# Row 9: Not found (lowercase), no G
$ws.Range("B9").Value = "Not found"
$ws.Range("C9").Value = "Not found"
$ws.Range("D9").Value = "Not found"
$ws.Range("E9").Value = "Not found"
$ws.Range("F9").Value = "Not found"

# Row 10: Not found (lowercase), no G
$ws.Range("B10").Value = "Not found"
$ws.Range("C10").Value = "Not found"
$ws.Range("D10").Value = "Not found"
$ws.Range("E10").Value = "Not found"
$ws.Range("F10").Value = "Not found"

# Row 11: Not found (lowercase), no G
$ws.Range("B11").Value = "Not found"
$ws.Range("C11").Value = "Not found"
$ws.Range("D11").Value = "Not found"
$ws.Range("E11").Value = "Not found"
$ws.Range("F11").Value = "Not found"
